$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old combined "Table_variability" test row (row 54) - rows below shift up.
$ws.Rows.Item(54).Delete()

# Append the two new, more specific variability-table test rows at the end (77, 78).
$ws.Cells.Item(77, 1).Value = "Table_variability1"
$ws.Cells.Item(77, 2).Value = "Test variability table for scale=252, digits=4"
$ws.Cells.Item(77, 3).Value = "table_variability_test1"

$ws.Cells.Item(78, 1).Value = "Table_variability2"
$ws.Cells.Item(78, 2).Value = "Test variability table for scale=1, digits=8"
$ws.Cells.Item(78, 3).Value = "table_variability_test2"

$ws.Application.ActiveWindow.ScrollRow = 61
$ws.Range("D82").Select()
